$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target roster data (player, position, team) for rows 2..19
$data = @(
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Jeremy Sochan", "SF,PF,C", "San Antonio Spurs"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves"),
    @("Quentin Grimes", "SG,SF", "Philadelphia 76ers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Amen Thompson", "PG,SG,SF,PF", "Houston Rockets"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
